# Reprocess the curated dimensions:
# The "aragon" column (D) used to be modeled as a curated SDMX dimension
# (iaest-dimension:aragon / skos:Concept / mapping-aragon.xlsx); it is now
# modeled directly as a reference area (sdmx-dimension:refArea / URI-Comunidad),
# matching how provincia-nombre (C) and comarca-nombre (G) are modeled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: iaest-dimension:aragon -> sdmx-dimension:refArea
$ws.Range("D2").Value = "sdmx-dimension:refArea"

# D4: skos:Concept -> URI-Comunidad
$ws.Range("D4").Value = "URI-Comunidad"

# D5: remove the mapping-aragon.xlsx cell entirely (no mapping file needed anymore)
$ws.Range("D5").Clear()
